$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 86 (existing rows 86-193 shift down to 87-194)
$ws.Rows("86:86").Insert()

# Populate the newly inserted row 86 with its data
$ws.Range("A86").Value = 4
$ws.Range("B86").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C86").Value = "Los Lagos"
$ws.Range("D86").Value = 44638
$ws.Range("E86").Value = 10
$ws.Range("F86").Value = 100112039
$ws.Range("G86").Value = "Ciboulette"
$ws.Range("H86").Value = "Sin especificar"
$ws.Range("I86").Value = "Primera"
$ws.Range("J86").Value = 120
$ws.Range("K86").Value = 6000
$ws.Range("L86").Value = 6000
$ws.Range("M86").Value = 6000
$ws.Range("N86").Value = '$/docena de atados'
$ws.Range("O86").Value = "Provincia de Cautín"
$ws.Range("P86").Value = 2000
$ws.Range("Q86").Value = 3
$ws.Range("R86").Value = "Hortaliza"
